$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 917
$ws.Range("F4").Value = 1083
$ws.Range("F5").Value = 1552
$ws.Range("F7").Value = 685
$ws.Range("F8").Value = 12362
$ws.Range("F10").Value = 2194
$ws.Range("F15").Value = 1237
$ws.Range("F17").Value = 280
$ws.Range("F19").Value = 681
$ws.Range("F21").Value = 2932
$ws.Range("F22").Value = 772
$ws.Range("F23").Value = 4260
$ws.Range("F24").Value = 4260
$ws.Range("F25").Value = 1140
$ws.Range("F26").Value = 872
$ws.Range("F30").Value = 1076
$ws.Range("F32").Value = 111
$ws.Range("F33").Value = 276
$ws.Range("F38").Value = 4470
$ws.Range("F40").Value = 4582
$ws.Range("F44").Value = 88
$ws.Range("F49").Value = 4114
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 66
$ws.Range("F5").Value = 102
$ws.Range("F7").Value = 44
$ws.Range("F9").Value = 3
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 472
$ws.Range("F4").Value = 90
$ws.Range("G4").Value = 39
$ws.Range("F5").Value = 12
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 472
$ws.Range("F4").Value = 90
$ws.Range("G4").Value = 39
$ws.Range("F6").Value = 917
$ws.Range("F7").Value = 1552
$ws.Range("F9").Value = 685
$ws.Range("F10").Value = 12362
$ws.Range("F11").Value = 2194
$ws.Range("F13").Value = 1237
$ws.Range("F14").Value = 280
$ws.Range("F15").Value = 681
$ws.Range("F17").Value = 2932
$ws.Range("F18").Value = 772
$ws.Range("F19").Value = 66
$ws.Range("F20").Value = 4260
$ws.Range("F21").Value = 1140
$ws.Range("F22").Value = 102
$ws.Range("F23").Value = 872
$ws.Range("F25").Value = 44
$ws.Range("F27").Value = 1076
$ws.Range("F29").Value = 111
$ws.Range("F30").Value = 276
$ws.Range("F33").Value = 4470
$ws.Range("F34").Value = 4582
$ws.Range("F38").Value = 88
$ws.Range("F45").Value = 4114
